$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels to reflect renamed columns:
#   I1: "Transaction Type" -> "Connection Type"
#   D1: "Marketplace ID" -> "Marketplace Name"
$ws.Range("I1").Value = "Connection Type"
$ws.Range("D1").Value = "Marketplace Name"

# Update the active selection to D1 (matches saved sheetView selection)
$ws.Range("D1").Select()
